$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy the existing row 15 (the old weekly record) down into new row 16,
# preserving its values/number formatting before row 15 gets overwritten with
# the new data.
$ws.Range("A16").Value2 = $ws.Range("A15").Value2
$ws.Range("B16").Value2 = $ws.Range("B15").Value2
$ws.Range("C16").Value2 = $ws.Range("C15").Value2
$ws.Range("D16").Value2 = $ws.Range("D15").Value2
$ws.Range("D16").NumberFormat = $ws.Range("D15").NumberFormat
$ws.Range("E16").Value2 = $ws.Range("E15").Value2
$ws.Range("F16").Value2 = $ws.Range("F15").Value2
$ws.Range("G16").Value2 = $ws.Range("G15").Value2
$ws.Range("H16").Value2 = $ws.Range("H15").Value2
$ws.Range("I16").Value2 = $ws.Range("I15").Value2
$ws.Range("J16").Value2 = $ws.Range("J15").Value2
$ws.Range("K16").Value2 = $ws.Range("K15").Value2
$ws.Range("L16").Value2 = $ws.Range("L15").Value2
$ws.Range("M16").Value2 = $ws.Range("M15").Value2
$ws.Range("N16").Value2 = $ws.Range("N15").Value2
$ws.Range("O16").Value2 = $ws.Range("O15").Value2
$ws.Range("P16").Value2 = $ws.Range("P15").Value2
$ws.Range("Q16").Value2 = $ws.Range("Q15").Value2
$ws.Range("R16").Value2 = $ws.Range("R15").Value2

# Step 2: update row 15 with the new weekly record's values.
$ws.Range("D15").Value2 = 44782
$ws.Range("J15").Value2 = 120
$ws.Range("K15").Value2 = 24000
$ws.Range("L15").Value2 = 24000
$ws.Range("M15").Value2 = 24000
$ws.Range("P15").Value2 = 1600
